# The workbook's data rows (2-17) got shuffled: each destination row now
# holds the data that used to live in a different source row (row 13 is
# the only one that stays put). Implement this as a pure row-content
# permutation over columns A:AY, preserving every field exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (i.e. new row N gets the data that used
# to be in row Map[N])
$map = @{
    2  = 3
    3  = 2
    4  = 17
    5  = 16
    6  = 5
    7  = 10
    8  = 11
    9  = 12
    10 = 6
    11 = 14
    12 = 4
    13 = 13
    14 = 15
    15 = 8
    16 = 7
    17 = 9
}

# Snapshot every source row's full A:AY contents before writing anything,
# so overlapping cycles in the permutation don't clobber data we still
# need to read.
$snapshot = @{}
for ($r = 2; $r -le 17; $r++) {
    $rng = $ws.Range("A" + $r + ":AY" + $r)
    $snapshot[$r] = $rng.Value()
}

# Columns holding a literal "YYYY-MM-DD" text value (Startdatum=Y=25,
# Slutdatum=AA=27). Excel's COM value-setter auto-detects those as real
# dates and silently coerces them to date serials; a leading apostrophe
# forces the literal-text interpretation instead, keeping the round trip
# lossless (no NumberFormat change either).
$dateTextCols = @(25, 27)

$colLetters = @(
    "A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z",
    "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY"
)

for ($destRow = 2; $destRow -le 17; $destRow++) {
    $srcRow = $map[$destRow]
    $vals = $snapshot[$srcRow]
    $guardedCols = @()

    for ($c = 1; $c -le 51; $c++) {
        $v = $vals[1, $c]
        if ($v -is [string]) {
            if ($dateTextCols -contains $c -and $v -match '^\d{4}-\d{2}-\d{2}$') {
                # guard literal dates from being reinterpreted as date serials
                $vals[1, $c] = "'" + $v
                $guardedCols += $c
            } elseif ($v -eq "") {
                # an empty string means the source cell genuinely exists
                # but is blank; plain "" would clear the cell entirely
                # (losing the distinction vs. a cell that never existed),
                # so write a bare apostrophe which round-trips to "".
                $vals[1, $c] = "'"
                $guardedCols += $c
            }
        }
    }

    $destRange = $ws.Range("A" + $destRow + ":AY" + $destRow)
    $destRange.Value = $vals

    # The leading apostrophe forces text entry the same way Excel's UI
    # quote-prefix does, which stamps the cell with a "quotePrefix" style.
    # Strip that back off (per guarded cell, not the whole row, so we
    # don't materialise brand-new empty cell nodes for columns that were
    # genuinely absent) so formatting stays identical to the source.
    foreach ($c in $guardedCols) {
        $ws.Range($colLetters[$c - 1] + $destRow).Style = "Normal"
    }
}
